# TC12_Canine_Filter_Breed-Boxer.xlsx — "corrected ICDC Breed 1-14 scripts"
#
# The FilesTab Cypher query (sheet "startup", cell B4) is corrected by
# dropping the two lines that return `File Type` and `Breed` (those
# properties are not meaningful/available for the Files result set), while
# every other line of the query is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFilesQuery = "MATCH (f:file)-->(parent)`n" +
                 "WITH DISTINCT f, parent`n" +
                 "MATCH (f)-[*]->(c:case)<--(demo:demographic)`n" +
                 "WHERE demo.breed IN ['Boxer'] `n" +
                 "OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`n" +
                 "OPTIONAL MATCH (samp:sample)-->(c)`n" +
                 "WITH DISTINCT f, parent, c, demo, diag, s`n" +
                 "RETURN  coalesce(f.file_name, '') AS ``File Name``,`n" +
                 "        coalesce(labels(parent)[0], '') AS ``Association``,`n" +
                 "        coalesce(f.file_description, '') AS ``Description``,`n" +
                 "        coalesce(f.file_format, '') AS ``Format``,`n" +
                 "        coalesce(f.file_size, '') AS ``Size``,`n" +
                 "        coalesce(c.case_id, '') AS ``Case ID``,`n" +
                 "        coalesce(diag.disease_term,'') AS Diagnosis , `n" +
                 "        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newFilesQuery

# The shorter query text reflows to fewer wrapped lines, so the row shrinks
# (246.5 -> 217.5) to match the new autofit content height.
$ws.Rows.Item(4).RowHeight = 217.5

# Selection moved from C3 to B4 (the edited cell) when the file was saved.
$ws.Activate() | Out-Null
$ws.Range("B4").Select() | Out-Null
